# Update the single data row (row 2) in the "UserDetails" sheet of TestData.xlsx
# with a new fake user record (Rupert Kling), replacing the previous
# "Crystal Zulauf" record. This mirrors the updated 'SelectDataInExcel'
# Apache POI based test-data generation described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Rupert"
$ws.Range("B2").Value = "Kling"
$ws.Range("C2").Value = "RupertKling30576"
$ws.Range("D2").Value = "vlcwwpr64ji2fa"
$ws.Range("E2").Value = "yelena.abshire@hotmail.com"
$ws.Range("F2").Value = "875-308-4141"
